$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row at position 13 (shifts rows 13-35 down to 14-36)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the UPDATE SCHEDULE command
$ws.Range("A13").Value = "UPDATE SCHEDULE schedName ‘schedule’"
$ws.Range("B13").Formula = "=LEFT(A13,SEARCH("" "",A13)-1)"
$ws.Range("C13").Value = "PUT"
$ws.Range("D13").Value = "schedules/{schedName}/body"

# Give the new row the same vertical-top alignment formatting used elsewhere
$ws.Range("A13,D13").VerticalAlignment = -4160

# Update the hidden AutoFilter defined name range to include the new row
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Commands!`$B`$1:`$B`$36"
